# Update the two key model inputs:
#  - Main!N3: shares outstanding assumption, 55 -> 48
#  - Model!V39: Q126 Revenue forecast input, 40000 -> 33000
# All other changed cells in the diff are formulas that depend on these
# two inputs and recalculate automatically.

$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("Main")
$wsModel = $wb.Worksheets.Item("Model")

$wsMain.Range("N3").Value = 48
$wsModel.Range("V39").Value = 33000

# Reproduce the final selection/active-sheet state recorded in the file:
# user last worked in Main!N3, then moved to the Model tab and selected Y53,
# leaving the Model sheet as the active tab.
$wsMain.Activate()
$null = $wsMain.Range("N3").Select()

$wsModel.Activate()
$null = $wsModel.Range("Y53").Select()
